$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(44326, 1, 8, 53.02578378736661),
    @(44327, 0, 7, 46.39756081394578),
    @(44328, 0, 6, 39.76933784052495),
    @(44329, 4, 10, 66.28222973420826)
)

$startRow = 252
$formatSourceRow = $startRow - 1

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]

    # Copy formatting (style) from the row above, like the existing rows use.
    $ws.Cells.Item($formatSourceRow, 1).Copy()
    $ws.Cells.Item($row, 1).PasteSpecial(-4122)

    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
    $ws.Cells.Item($row, 4).Value = $vals[3]
}

$excel.CutCopyMode = $false
